$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.735.71"
$ws.Range("E2").Value = "  -1.73%  "
$ws.Range("D3").Value = "1.758.47"
$ws.Range("E3").Value = "  -2.13%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "326.84"
$ws.Range("E5").Value = "  -2.49%  "
$ws.Range("E6").Value = "  -0.05%  "
$ws.Range("D7").Value = "0.4436"
$ws.Range("E7").Value = "  -1.22%  "
$ws.Range("D8").Value = "0.3751"
$ws.Range("E8").Value = "  +1.00%  "
$ws.Range("E9").Value = "  +0.69%  "
$ws.Range("D10").Value = "0.07652"
$ws.Range("E10").Value = "  +1.10%  "
$ws.Range("E11").Value = "  -1.48%  "
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("E13").Value = "  -2.89%  "
$ws.Range("D14").Value = "6.203"
$ws.Range("E14").Value = "  -1.46%  "
$ws.Range("D15").Value = "7.453"
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("D16").Value = "1.757.23"
$ws.Range("E16").Value = "  -1.98%  "
$ws.Range("D17").Value = "0.00001076"
$ws.Range("E17").Value = "  -1.21%  "
$ws.Range("D18").Value = "89.00"
$ws.Range("E18").Value = "  +9.76%  "
$ws.Range("D19").Value = "0.06216"
$ws.Range("E19").Value = "  -7.89%  "
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("E21").Value = "  -0.79%  "
$ws.Range("E22").Value = "  -2.91%  "
$ws.Range("D23").Value = "0.5323"
$ws.Range("E23").Value = "  -3.15%  "
$ws.Range("D24").Value = "27.763.01"
$ws.Range("E24").Value = "  -1.58%  "
$ws.Range("D25").Value = "11.63"
$ws.Range("E25").Value = "  -1.41%  "
$ws.Range("D26").Value = "2.315"
$ws.Range("E26").Value = "  -4.17%  "
$ws.Range("D27").Value = "20.75"
$ws.Range("E27").Value = "  +1.06%  "
$ws.Range("D28").Value = "153.83"
$ws.Range("E28").Value = "  +1.40%  "
$ws.Range("D29").Value = "2.362"
$ws.Range("D30").Value = "1.955.63"
$ws.Range("E30").Value = "  -2.11%  "
$ws.Range("D31").Value = "128.33"
$ws.Range("E31").Value = "  -3.46%  "
$ws.Range("E32").Value = "  -1.76%  "
$ws.Range("D33").Value = "0.09362"
$ws.Range("E33").Value = "  -0.41%  "
$ws.Range("D34").Value = "5.767"
$ws.Range("E34").Value = "  -0.56%  "
$ws.Range("D35").Value = "3.659"
$ws.Range("E35").Value = "  -9.45%  "
$ws.Range("D36").Value = "12.72"
$ws.Range("E36").Value = "  +5.52%  "
$ws.Range("D37").Value = "0.2186"
$ws.Range("E37").Value = "  -7.25%  "
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("D39").Value = "0.06155"
$ws.Range("E39").Value = "  -2.50%  "
$ws.Range("D42").Value = "1.204"
$ws.Range("E42").Value = "  -0.42%  "
$ws.Range("D43").Value = "8.011"
$ws.Range("E43").Value = "  -3.90%  "
$ws.Range("E44").Value = "  -4.33%  "
$ws.Range("D45").Value = "1.001"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").Value = "13.79"
$ws.Range("E46").Value = "  -3.44%  "
$ws.Range("D47").Value = "0.6022"
$ws.Range("E47").Value = "  -0.92%  "
$ws.Range("D48").Value = "3.767"
$ws.Range("E48").Value = "  -1.64%  "
$ws.Range("D49").Value = "126.38"
$ws.Range("E49").Value = "  -2.80%  "
$ws.Range("D50").Value = "1.999"
$ws.Range("E50").Value = "  -1.40%  "

# Row 40 and 41 swap (InternetComputer(DFINITY) <-> TheSandbox) with updated values
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "0.6499"
$ws.Range("E40").Value = "  -0.94%  "

$ws.Range("B41").Value = "InternetComputer(DFINITY)"
$ws.Range("C41").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D41").Value = "5.095"
$ws.Range("E41").Value = "  -2.12%  "

# Row 51: Cronos replaced with EOS
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").Value = "1.138"
$ws.Range("E51").Value = "  -1.81%  "
